$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F29").Value = 12
$ws.Range("G29").Value = 368.88
$ws.Range("B34").Value = 10364.71
$ws.Range("F38").Value = 58
$ws.Range("G38").Value = 1811.34
$ws.Range("F43").Value = 13
$ws.Range("G43").Value = 2507.57
$ws.Range("F44").Value = 34
$ws.Range("G44").Value = 1199.86
$ws.Range("B63").Value = 37503.94
$ws.Range("F81").Value = 55
$ws.Range("G81").Value = 2189
$ws.Range("B82").Value = 24185.08
$ws.Range("F92").Value = 214
$ws.Range("G92").Value = 17370.38
$ws.Range("F100").Value = 278
$ws.Range("G100").Value = 17708.6
$ws.Range("F111").Value = 136
$ws.Range("G111").Value = 13934.56
$ws.Range("F113").Value = 0
$ws.Range("G113").Value = 0
$ws.Range("F114").Value = 14
$ws.Range("G114").Value = 707.7
$ws.Range("F120").Value = 80
$ws.Range("G120").Value = 3784
$ws.Range("F121").Value = 84
$ws.Range("G121").Value = 1634.64
$ws.Range("B122").Value = 236324.79
$ws.Range("F157").Value = 135
$ws.Range("G157").Value = 6010.2
$ws.Range("B160").Value = 21663.93
$ws.Range("B165").Value = 57756
$ws.Range("B166").Value = 53925
$ws.Range("F190").Value = 16
$ws.Range("G190").Value = 1336
$ws.Range("F195").Value = 25
$ws.Range("G195").Value = 2050
$ws.Range("F205").Value = 82
$ws.Range("G205").Value = 3276.72
$ws.Range("B206").Value = 15637.87
$ws.Range("F210").Value = 129
$ws.Range("G210").Value = 4002.87
$ws.Range("B218").Value = 13611.65
$ws.Range("F235").Value = 20
$ws.Range("G235").Value = 1696.6
$ws.Range("B238").Value = 8334.43
$ws.Range("F270").Value = 71
$ws.Range("G270").Value = 5725.44
$ws.Range("B300").Value = 115286.2
$ws.Range("B304").Value = 61610
$ws.Range("D304").Value = 102.71
$ws.Range("E304").Value = 122.71
$ws.Range("F304").Value = 176
$ws.Range("G304").Value = 18076.96
$ws.Range("B305").Value = 57077
$ws.Range("D305").Value = 93.08
$ws.Range("E305").Value = 111.2
$ws.Range("F305").Value = 1
$ws.Range("G305").Value = 93.08
$ws.Range("F307").Value = 73
$ws.Range("G307").Value = 7817.57
$ws.Range("F314").Value = 182
$ws.Range("G314").Value = 20789.86
$ws.Range("F317").Value = 86
$ws.Range("G317").Value = 11788.02
$ws.Range("F333").Value = 146
$ws.Range("G333").Value = 16346.16
$ws.Range("F341").Value = 15
$ws.Range("G341").Value = 1332
$ws.Range("F344").Value = 113
$ws.Range("G344").Value = 11424.3
$ws.Range("F351").Value = 57
$ws.Range("G351").Value = 3370.41
$ws.Range("F359").Value = 8
$ws.Range("G359").Value = 4207.6
$ws.Range("F362").Value = 190
$ws.Range("G362").Value = 13358.9
$ws.Range("F366").Value = 34
$ws.Range("G366").Value = 2095.76
$ws.Range("B370").Value = 335877.41
$ws.Range("F412").Value = 32
$ws.Range("G412").Value = 31033.92
$ws.Range("B413").Value = 31033.92
$ws.Range("B417").Value = 58047
$ws.Range("D417").Value = 105.54
$ws.Range("E417").Value = 126.1
$ws.Range("F417").Value = 62
$ws.Range("G417").Value = 6543.48
$ws.Range("B418").Value = 47097
$ws.Range("D418").Value = 112.28
$ws.Range("E418").Value = 134.16
$ws.Range("F418").Value = 15
$ws.Range("G418").Value = 1684.2
$ws.Range("F423").Value = 20
$ws.Range("G423").Value = 1932
$ws.Range("F424").Value = 51
$ws.Range("G424").Value = 1897.71
$ws.Range("F425").Value = 58
$ws.Range("G425").Value = 1426.22
$ws.Range("B428").Value = 44993.83
$ws.Range("F431").Value = 16
$ws.Range("G431").Value = 798.24
$ws.Range("F432").Value = 60
$ws.Range("G432").Value = 3109.2
$ws.Range("F433").Value = 49
$ws.Range("G433").Value = 4814.25
$ws.Range("F434").Value = 34
$ws.Range("G434").Value = 1378.36
$ws.Range("F437").Value = 42
$ws.Range("G437").Value = 1571.64
$ws.Range("F439").Value = 86
$ws.Range("G439").Value = 2859.5
$ws.Range("F440").Value = 55
$ws.Range("G440").Value = 3324.75
$ws.Range("B445").Value = 40186.17
$ws.Range("F466").Value = 598
$ws.Range("G466").Value = 8043.1
$ws.Range("F468").Value = 603
$ws.Range("G468").Value = 7724.43
$ws.Range("F472").Value = 303
$ws.Range("G472").Value = 3881.43
$ws.Range("F473").Value = 374
$ws.Range("G473").Value = 7379.02
$ws.Range("F474").Value = 423
$ws.Range("G474").Value = 2783.34
$ws.Range("F483").Value = 611
$ws.Range("G483").Value = 9000.030000000001
$ws.Range("B484").Value = 104103.68
$ws.Range("F486").Value = 49
$ws.Range("G486").Value = 1811.53
$ws.Range("F488").Value = 1
$ws.Range("G488").Value = 217.08
$ws.Range("F490").Value = 109
$ws.Range("G490").Value = 3794.29
$ws.Range("B491").Value = 8822.91
$ws.Range("F502").Value = 28
$ws.Range("G502").Value = 686
$ws.Range("B509").Value = 38158.42
$ws.Range("F547").Value = 29
$ws.Range("G547").Value = 1513.22
$ws.Range("B556").Value = 26016.98
$ws.Range("F593").Value = 3
$ws.Range("G593").Value = 852.3
$ws.Range("F601").Value = 23
$ws.Range("G601").Value = 1909.46
$ws.Range("B602").Value = 36141.97
$ws.Range("F607").Value = 102
$ws.Range("G607").Value = 9684.9
$ws.Range("F608").Value = 69
$ws.Range("G608").Value = 1876.8
$ws.Range("B611").Value = 61730.42
$ws.Range("F618").Value = 0
$ws.Range("G618").Value = 0
$ws.Range("B628").Value = 63124.65
$ws.Range("F630").Value = 44
$ws.Range("G630").Value = 1456.84
$ws.Range("F631").Value = 200
$ws.Range("G631").Value = 3172
$ws.Range("F632").Value = 97
$ws.Range("G632").Value = 3211.67
$ws.Range("F634").Value = 0
$ws.Range("G634").Value = 0
$ws.Range("F635").Value = 84
$ws.Range("G635").Value = 3627.12
$ws.Range("F636").Value = 68
$ws.Range("G636").Value = 2251.48
$ws.Range("B639").Value = 28443.49
$ws.Range("F653").Value = 8
$ws.Range("G653").Value = 4618.32
$ws.Range("B658").Value = 9121.65
$ws.Range("F704").Value = 5
$ws.Range("G704").Value = 553.15
$ws.Range("F709").Value = 117
$ws.Range("G709").Value = 10609.56
$ws.Range("F711").Value = 39
$ws.Range("G711").Value = 1616.55
$ws.Range("B719").Value = 479624.68
$ws.Range("F728").Value = 224
$ws.Range("G728").Value = 6746.88
$ws.Range("F730").Value = 27
$ws.Range("G730").Value = 2778.57
$ws.Range("B737").Value = 18322.75
$ws.Range("F742").Value = 75
$ws.Range("G742").Value = 2805
$ws.Range("B745").Value = 7256.18
$ws.Range("F790").Value = 771
$ws.Range("G790").Value = 125757.81
$ws.Range("F792").Value = 207
$ws.Range("G792").Value = 15967.98
$ws.Range("F793").Value = 28
$ws.Range("G793").Value = 4140.08
$ws.Range("F794").Value = 63
$ws.Range("G794").Value = 4252.5
$ws.Range("B795").Value = 150684.11
$ws.Range("B801").Value = 2960027
$ws.Range("B802").Value = 2960027
